{"js": "// Replace the 100 arithmetic-expression cell values in the single table\n// with their updated counterparts, preserving all run/paragraph formatting.\nconst newValues = [\n  [\"17+35=\", \"52+28=\", \"51+32=\", \"37+55=\", \"44+45=\"],\n  [\"28-15=\", \"42+46=\", \"31+9=\", \"81-42=\", \"65-4=\"],\n  [\"15+4=\", \"83+14=\", \"9+55=\", \"66-1=\", \"63-10=\"],\n  [\"25+15=\", \"24+72=\", \"73-29=\", \"82-11=\", \"44-9=\"],\n  [\"14+53=\", \"44+43=\", \"10+81=\", \"34+19=\", \"66-8=\"],\n  [\"79-28=\", \"13-5=\", \"4+1=\", \"44-35=\", \"95-41=\"],\n  [\"92-58=\", \"35+12=\", \"61-42=\", \"23+52=\", \"62+34=\"],\n  [\"74-58=\", \"17+10=\", \"71-55=\", \"30-9=\", \"89-77=\"],\n  [\"27-21=\", \"13+33=\", \"44-24=\", \"67-1=\", \"52+44=\"],\n  [\"74+15=\", \"49-18=\", \"77-46=\", \"65-45=\", \"51+43=\"],\n  [\"75-34=\", \"86-73=\", \"14+19=\", \"11+16=\", \"24+45=\"],\n  [\"89-2=\", \"9+33=\", \"99-17=\", \"31+40=\", \"20+11=\"],\n  [\"69-0=\", \"22+39=\", \"26+51=\", \"9+69=\", \"27+27=\"],\n  [\"79-51=\", \"62-38=\", \"86+11=\", \"50-17=\", \"65-18=\"],\n  [\"28+58=\", \"23+10=\", \"90-50=\", \"98-83=\", \"60-51=\"],\n  [\"20+0=\", \"44+27=\", \"20+0=\", \"69-58=\", \"97-60=\"],\n  [\"24+30=\", \"9+73=\", \"97-94=\", \"72+17=\", \"81-56=\"],\n  [\"43+49=\", \"22-0=\", \"26+31=\", \"71+12=\", \"75-37=\"],\n  [\"28-13=\", \"6+40=\", \"10+86=\", \"85+11=\", \"25+53=\"],\n  [\"40+49=\", \"55-17=\", \"50+23=\", \"60+35=\", \"67-24=\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nfor (let r = 0; r < newValues.length; r++) {\n  for (let c = 0; c < newValues[r].length; c++) {\n    const cell = table.getCell(r, c);\n    cell.body.paragraphs.load(\"items\");\n    await context.sync();\n\n    const paragraph = cell.body.paragraphs.items[0];\n    const range = paragraph.getRange();\n    range.insertText(newValues[r][c], Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the 100 arithmetic-expression cell values in the single table\n# with their updated counterparts, preserving all run/paragraph formatting.\n$newValues = @(\n  @(\"17+35=\", \"52+28=\", \"51+32=\", \"37+55=\", \"44+45=\"),\n  @(\"28-15=\", \"42+46=\", \"31+9=\", \"81-42=\", \"65-4=\"),\n  @(\"15+4=\", \"83+14=\", \"9+55=\", \"66-1=\", \"63-10=\"),\n  @(\"25+15=\", \"24+72=\", \"73-29=\", \"82-11=\", \"44-9=\"),\n  @(\"14+53=\", \"44+43=\", \"10+81=\", \"34+19=\", \"66-8=\"),\n  @(\"79-28=\", \"13-5=\", \"4+1=\", \"44-35=\", \"95-41=\"),\n  @(\"92-58=\", \"35+12=\", \"61-42=\", \"23+52=\", \"62+34=\"),\n  @(\"74-58=\", \"17+10=\", \"71-55=\", \"30-9=\", \"89-77=\"),\n  @(\"27-21=\", \"13+33=\", \"44-24=\", \"67-1=\", \"52+44=\"),\n  @(\"74+15=\", \"49-18=\", \"77-46=\", \"65-45=\", \"51+43=\"),\n  @(\"75-34=\", \"86-73=\", \"14+19=\", \"11+16=\", \"24+45=\"),\n  @(\"89-2=\", \"9+33=\", \"99-17=\", \"31+40=\", \"20+11=\"),\n  @(\"69-0=\", \"22+39=\", \"26+51=\", \"9+69=\", \"27+27=\"),\n  @(\"79-51=\", \"62-38=\", \"86+11=\", \"50-17=\", \"65-18=\"),\n  @(\"28+58=\", \"23+10=\", \"90-50=\", \"98-83=\", \"60-51=\"),\n  @(\"20+0=\", \"44+27=\", \"20+0=\", \"69-58=\", \"97-60=\"),\n  @(\"24+30=\", \"9+73=\", \"97-94=\", \"72+17=\", \"81-56=\"),\n  @(\"43+49=\", \"22-0=\", \"26+31=\", \"71+12=\", \"75-37=\"),\n  @(\"28-13=\", \"6+40=\", \"10+86=\", \"85+11=\", \"25+53=\"),\n  @(\"40+49=\", \"55-17=\", \"50+23=\", \"60+35=\", \"67-24=\"),\n)\n\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\nfor ($r = 1; $r -le $newValues.Length; $r++) {\n  $row = $newValues[$r - 1]\n  for ($c = 1; $c -le $row.Length; $c++) {\n    $table.Cell($r, $c).Range.Text = $row[$c - 1]\n  }\n}\n"}
